# Update cryptos list (prices & 1h volume %) per Wed May 15 13:41:37 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.434.45'
$ws.Range("E2").Value = '  +3.98%  '

# Row 3
$ws.Range("D3").Value = '2.974.09'
$ws.Range("E3").Value = '  +2.24%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").Value = '''581.33'
$ws.Range("E5").Value = '  +0.36%  '

# Row 6
$ws.Range("D6").Value = '''152.31'
$ws.Range("E6").Value = '  +4.18%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").Value = '2.973.80'
$ws.Range("E8").Value = '  +2.24%  '

# Row 9
$ws.Range("D9").Value = '''0.510'
$ws.Range("E9").Value = '  +0.85%  '

# Row 10
$ws.Range("D10").Value = '''6.99'
$ws.Range("E10").Value = '  +3.96%  '

# Row 11
$ws.Range("E11").Value = '  -0.12%  '

# Row 12
$ws.Range("D12").Value = '''0.447'
$ws.Range("E12").Value = '  +2.99%  '

# Row 13
$ws.Range("D13").Value = '''0.0000241'
$ws.Range("E13").Value = '  +1.70%  '

# Row 14
$ws.Range("D14").Value = '''34.37'
$ws.Range("E14").Value = '  +5.02%  '

# Row 15
$ws.Range("E15").Value = '  +0.72%  '

# Row 16
$ws.Range("D16").Value = '3.463.17'
$ws.Range("E16").Value = '  +2.05%  '

# Row 17
$ws.Range("D17").Value = '64.261.96'
$ws.Range("E17").Value = '  +3.71%  '

# Row 18
$ws.Range("D18").Value = '''6.90'
$ws.Range("E18").Value = '  +3.64%  '

# Row 19
$ws.Range("D19").Value = '2.965.17'
$ws.Range("E19").Value = '  +1.96%  '

# Row 20
$ws.Range("D20").Value = '''456.01'
$ws.Range("E20").Value = '  +4.67%  '

# Row 21
$ws.Range("D21").Value = '''13.64'
$ws.Range("E21").Value = '  +2.30%  '

# Row 22
$ws.Range("D22").Value = '''0.675'
$ws.Range("E22").Value = '  +2.24%  '

# Row 23
$ws.Range("D23").Value = '''7.15'
$ws.Range("E23").Value = '  +2.78%  '

# Row 24
$ws.Range("D24").Value = '''80.51'
$ws.Range("E24").Value = '  +0.68%  '

# Row 25
$ws.Range("D25").Value = '''10.99'
$ws.Range("E25").Value = '  +7.83%  '

# Row 26
$ws.Range("D26").Value = '''12.26'
$ws.Range("E26").Value = '  +2.67%  '

# Row 27
$ws.Range("D27").Value = '''2.20'
$ws.Range("E27").Value = '  +6.90%  '

# Row 28
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
$ws.Range("D29").Value = '''7.63'
$ws.Range("E29").Value = '  +7.71%  '

# Row 30
$ws.Range("E30").Value = '  -2.26%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''2.14'
$ws.Range("E31").Value = '  +2.01%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.56'
$ws.Range("E32").Value = '  +0.52%  '

# Row 33
$ws.Range("D33").Value = '''0.110'
$ws.Range("E33").Value = '  +2.39%  '

# Row 34
$ws.Range("D34").Value = '''26.58'
$ws.Range("E34").Value = '  +3.00%  '

# Row 35
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  -0.14%  '

# Row 36
$ws.Range("D36").Value = '''0.973'
$ws.Range("E36").Value = '  +0.35%  '

# Row 37
$ws.Range("D37").Value = '''2.14'
$ws.Range("E37").Value = '  +8.29%  '

# Row 38
$ws.Range("D38").Value = '''5.63'
$ws.Range("E38").Value = '  +2.47%  '

# Row 39
$ws.Range("D39").Value = '''3.04'
$ws.Range("E39").Value = '  -1.38%  '

# Row 40
$ws.Range("D40").Value = '''49.08'
$ws.Range("E40").Value = '  -0.27%  '

# Row 41
$ws.Range("D41").Value = '''44.31'
$ws.Range("E41").Value = '  +15.10%  '

# Row 42
$ws.Range("D42").Value = '''0.119'
$ws.Range("E42").Value = '  +1.90%  '

# Row 43
$ws.Range("D43").Value = '''0.292'
$ws.Range("E43").Value = '  +8.03%  '

# Row 44
$ws.Range("D44").Value = '''8.32'
$ws.Range("E44").Value = '  +0.08%  '

# Row 45
$ws.Range("D45").Value = '''379.49'
$ws.Range("E45").Value = '  +10.33%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0351'
$ws.Range("E46").Value = '  +4.24%  '

# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.749.95'
$ws.Range("E47").Value = '  +2.28%  '

# Row 48
$ws.Range("D48").Value = '''134.58'
$ws.Range("E48").Value = '  -0.21%  '

# Row 49
$ws.Range("E49").Value = '  -0.01%  '

# Row 50
$ws.Range("E50").Value = '  +1.90%  '

# Row 51
$ws.Range("D51").Value = '''0.000217'
$ws.Range("E51").Value = '  +6.62%  '
